$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1275099798234569
$ws.Range("C2").Value = 1.51412274867419
$ws.Range("D2").Value = 9.324399010553055
$ws.Range("E2").Value = 3.053587891407918
$ws.Range("F2").Value = 3.122720766858579
$ws.Range("G2").Value = 22
$ws.Range("B3").Value = 0.02989055876380673
$ws.Range("C3").Value = 1.469662598613318
$ws.Range("D3").Value = 8.914287453216375
$ws.Range("E3").Value = 2.985680400380519
$ws.Range("F3").Value = 3.059258686037979
$ws.Range("G3").Value = 21
$ws.Range("B4").Value = -0.4518499121642435
$ws.Range("C4").Value = 0.921173439906215
$ws.Range("D4").Value = 3.928571563801758
$ws.Range("E4").Value = 1.982062452043769
$ws.Range("F4").Value = 1.980006172419652
$ws.Range("G4").Value = 20
$ws.Range("B5").Value = -0.02887301196522735
$ws.Range("C5").Value = 0.640656743365336
$ws.Range("D5").Value = 0.9487815114782999
$ws.Range("E5").Value = 0.9740541624972915
$ws.Range("F5").Value = 1.000305768600247
$ws.Range("G5").Value = 19
$ws.Range("B6").Value = -0.01274630076391603
$ws.Range("C6").Value = 0.6862276495140773
$ws.Range("D6").Value = 0.9980462345064838
$ws.Range("E6").Value = 0.9990226396366019
$ws.Range("F6").Value = 1.027902140965356
$ws.Range("G6").Value = 18
$ws.Range("B7").Value = -0.1249830523586888
$ws.Range("C7").Value = 0.5228446924603924
$ws.Range("D7").Value = 0.510271894425645
$ws.Range("E7").Value = 0.7143331816636023
$ws.Range("F7").Value = 0.7249598793997479
$ws.Range("G7").Value = 17
$ws.Range("B8").Value = -0.03908073210606425
$ws.Range("C8").Value = 0.476833024709525
$ws.Range("D8").Value = 0.409094268446219
$ws.Range("E8").Value = 0.6396047751902881
$ws.Range("F8").Value = 0.6593467189670075
$ws.Range("G8").Value = 16
$ws.Range("B9").Value = 0.04629975205653263
$ws.Range("C9").Value = 0.4102896313898315
$ws.Range("D9").Value = 0.3245141161252205
$ws.Range("E9").Value = 0.5696614048057148
$ws.Range("F9").Value = 0.5877047811049635
$ws.Range("G9").Value = 15
$ws.Range("B10").Value = 0.007747740397374882
$ws.Range("C10").Value = 0.3348337809727632
$ws.Range("D10").Value = 0.261718090368209
$ws.Range("E10").Value = 0.5115839035468268
$ws.Range("F10").Value = 0.5308348200579344
$ws.Range("G10").Value = 14
$ws.Range("B11").Value = 0.008384872907047361
$ws.Range("C11").Value = 0.3303106150213118
$ws.Range("D11").Value = 0.2394032953374836
$ws.Range("E11").Value = 0.4892885603991611
$ws.Range("F11").Value = 0.509192895683094
